# Applies the "Results now are ordered in correct levels" edit:
# The season factor levels are reordered (Winter, Spring, Summer instead of
# Spring, Summer, Winter), which shifts which term label/row each estimate
# row corresponds to, and the model's fixed-effect estimates are refreshed
# to match the corrected level ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; Term='(Intercept)'; E=[double]"-1.85411499144923"; F=[double]"0.0567349096960273"; G=[double]"-32.6803197781252"; H=[double]"2.97356868114431e-234" }
    @{ Row=3; Term='habitat_typeExposed/Low SAV'; E=[double]"0.292041300078359"; F=[double]"0.123039009651358"; G=[double]"2.37356673225739"; H=[double]"0.0176172092768298" }
    @{ Row=4; Term='habitat_typeMod/Dense SAV'; E=[double]"-0.0861170786027511"; F=[double]"0.0169913130950191"; G=[double]"-5.0683003792094"; H=[double]"0.000000401383631927643" }
    @{ Row=5; Term='habitat_typeShallow/Dense SAV'; E=[double]"0.0946530536493692"; F=[double]"0.0794277394570933"; G=[double]"1.1916876181589"; H=[double]"0.233383754834432" }
    @{ Row=6; Term='habitat_typeShallow/Low SAV'; E=[double]"0.0460185247374881"; F=[double]"0.01415972698173"; G=[double]"3.24995847708539"; H=[double]"0.00115421860013021" }
    @{ Row=7; Term='seasonWinter'; E=[double]"-0.228764292516279"; F=[double]"0.00895236443451118"; G=[double]"-25.553505354898"; H=[double]"0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000050195684143135" }
    @{ Row=8; Term='seasonSpring'; E=[double]"0.164500488019238"; F=[double]"0.0135562757853923"; G=[double]"12.1346371690444"; H=[double]"0.000000000000000000000000000000000692161188773339" }
    @{ Row=9; Term='seasonSummer'; E=[double]"0.255197283194904"; F=[double]"0.0136955109831708"; G=[double]"18.6336445210765"; H=[double]"0.0000000000000000000000000000000000000000000000000000000000000000000000000000171442597565997" }
    @{ Row=10; Term='habitat_typeExposed/Low SAV:seasonWinter'; E=[double]"0.0851206424651651"; F=[double]"0.131336357320149"; G=[double]"0.648111796322116"; H=[double]"0.516912644860807" }
    @{ Row=11; Term='habitat_typeMod/Dense SAV:seasonWinter'; E=[double]"0.0915433448082384"; F=[double]"0.0194387474665106"; G=[double]"4.70932321981912"; H=[double]"0.00000248540671398563" }
    @{ Row=12; Term='habitat_typeShallow/Dense SAV:seasonWinter'; E=[double]"0.521227784984202"; F=[double]"0.0816279014577631"; G=[double]"6.38541204264454"; H=[double]"0.000000000170936572716764" }
    @{ Row=13; Term='habitat_typeShallow/Low SAV:seasonWinter'; E=[double]"-0.0749208696098756"; F=[double]"0.016383101579395"; G=[double]"-4.57305774775293"; H=[double]"0.00000480657594265583" }
    @{ Row=14; Term='habitat_typeExposed/Low SAV:seasonSpring'; E=[double]"-0.0436976512199168"; F=[double]"0.129547596036978"; G=[double]"-0.337309626397419"; H=[double]"0.735883503687849" }
    @{ Row=15; Term='habitat_typeMod/Dense SAV:seasonSpring'; E=[double]"0.129527806961306"; F=[double]"0.0264026767632323"; G=[double]"4.90585890676368"; H=[double]"0.000000930193764502628" }
    @{ Row=16; Term='habitat_typeShallow/Dense SAV:seasonSpring'; E=[double]"-0.411190914942214"; F=[double]"0.081410665072834"; G=[double]"-5.05082367985009"; H=[double]"0.000000439908994957145" }
    @{ Row=17; Term='habitat_typeShallow/Low SAV:seasonSpring'; E=[double]"0.21488358836013"; F=[double]"0.0230240720598793"; G=[double]"9.33299669151817"; H=[double]"0.0000000000000000000102919046826049" }
    @{ Row=18; Term='habitat_typeExposed/Low SAV:seasonSummer'; E=[double]"0.0443660577753137"; F=[double]"0.137030642270912"; G=[double]"0.32376742194349"; H=[double]"0.746114122291077" }
    @{ Row=19; Term='habitat_typeMod/Dense SAV:seasonSummer'; E=[double]"-0.0778877183512656"; F=[double]"0.0250979892837383"; G=[double]"-3.10334495208871"; H=[double]"0.00191346467935494" }
    @{ Row=20; Term='habitat_typeShallow/Dense SAV:seasonSummer'; E=[double]"-0.452402862014736"; F=[double]"0.0820430264217772"; G=[double]"-5.51421469618839"; H=[double]"0.0000000350340719496958" }
    @{ Row=21; Term='habitat_typeShallow/Low SAV:seasonSummer'; E=[double]"-0.0965099941672309"; F=[double]"0.0191725825555898"; G=[double]"-5.03375035091937"; H=[double]"0.000000480975967108472" }
    @{ Row=22; Term='sd__(Intercept)'; E=[double]"0.250378129464447"; F=$null; G=$null; H=$null }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 4).Value = $row.Term
    $ws.Cells.Item($row.Row, 5).Value = $row.E
    if ($null -ne $row.F) { $ws.Cells.Item($row.Row, 6).Value = $row.F }
    if ($null -ne $row.G) { $ws.Cells.Item($row.Row, 7).Value = $row.G }
    if ($null -ne $row.H) { $ws.Cells.Item($row.Row, 8).Value = $row.H }
}
